# Applies the "3/14/2016" log-entry update:
#  - marks a few camelCase identifiers with spell/grammar proofing runs
#    (matching Word's automatic proofErr bookkeeping after an edit)
#  - appends the new 3/14/2016 log entries after the SOLUTION bullet
#  - relocates the _GoBack bookmark to the new last (empty) paragraph,
#    which is what Word does after the final edit in the document

function Insert-BodyXml($range, [string]$bodyXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

$d = $word.ActiveDocument

# --- Paragraph 2: "3/12/16 (6 hr)" -> flag "hr" as a spelling run ---
$p2 = $d.Paragraphs.Item(2)
$body2 = '<w:body><w:p>' +
  '<w:r><w:t xml:space="preserve">3/12/16 (6 </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>hr</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>)</w:t></w:r>' +
  '</w:p></w:body>'
Insert-BodyXml $p2.Range $body2

# --- Paragraph 4: "In StartActivity class, ..." -> flag "StartActivity" ---
$p4 = $d.Paragraphs.Item(4)
$body4 = '<w:body><w:p>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">In </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>StartActivity</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> class, there are two options since we have two buttons. Start new game and load a game. If new game button is pressed, it takes the existing format from our RAW resource. Then it is rendered and the data is extracted from that file.</w:t></w:r>' +
  '</w:p></w:body>'
Insert-BodyXml $p4.Range $body4

# --- Paragraph 11: "Changed the function readData(rawFile) to readData(rawFile, filename, indicator) ..." ---
$p11 = $d.Paragraphs.Item(11)
$body11 = '<w:body><w:p>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Changed the function </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>readData</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>(</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>rawFile</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve">) to </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>readData</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>(</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>rawFile</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>, filename, indicator) where indicator indicates whether to read from raw data or the filename the user provided.</w:t></w:r>' +
  '</w:p></w:body>'
Insert-BodyXml $p11.Range $body11

# --- Paragraph 12: "Checked in StartPageActivity ... fileAccess." ---
$p12 = $d.Paragraphs.Item(12)
$body12 = '<w:body><w:p>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Checked in </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>StartPageActivity</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> if the user entered filename already exists to open. If not, we tell the user about it and do not go through the entire trouble of reaching up to </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>fileAccess</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>.</w:t></w:r>' +
  '</w:p></w:body>'
Insert-BodyXml $p12.Range $body12

# --- Paragraph 14 (last paragraph, "SOLUTION: ...") through end of document:
#     keep its text as-is, then append the new 3/14/2016 entries, then
#     finish with an empty paragraph holding the relocated _GoBack bookmark ---
$p14 = $d.Paragraphs.Item(14)
$tail = $d.Range($p14.Range.Start, $d.Content.End)

$body14 = '<w:body>' +
  '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">SOLUTION: While writing, open the board and go through each lines. For each column, just add a space and append. For each row, add \n so that we can distinguish later. </w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:r><w:t>3/14/2016</w:t></w:r></w:p>' +
  '<w:p><w:r><w:tab/><w:t>(total &#8211; 2 hours)</w:t></w:r></w:p>' +
  '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Changed the </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>AlertDialog</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> type for handling the actual game of heads or tails. Then, once the user chooses it, there will be another dialog box that will give the result. Then, it is directed to the main game. </w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">In </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>MainActivity</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, things are changed now. I was under the impression that the new game had to be loaded from the RAW but it wasn&#8217;t. New game was completely new. So, I just separated the </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>fileAccess</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> part for the load game section.</w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Made Save Game button which saves the game as savedGame.txt file.</w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Previously I was using internal storage. So, changed it to external storage and formatted it with the given serialization.</w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>So, Lots of previous techniques were changed here.</w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">TO DO: Show the player&#8217;s next tile. Also option to see further in the stock. </w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>' +
  '</w:body>'
Insert-BodyXml $tail $body14

Write-Output "Edit complete."
